# Updates cryptos list values (price & volume change) per upstream data refresh,
# including the Uniswap/ShibaInu row-order swap (rows 19-20).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''41.069.52'
$ws.Range('E2').Value = '  -1.33%  '
$ws.Range('D3').Value = '''2.425.92'
$ws.Range('E3').Value = '  -1.87%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = '''317.67'
$ws.Range('E5').Value = '  -0.02%  '
$ws.Range('D6').Value = '''89.34'
$ws.Range('E6').Value = '  -3.28%  '
$ws.Range('D7').Value = '''0.538'
$ws.Range('E7').Value = '  -2.54%  '
$ws.Range('E9').Value = '  -3.29%  '
$ws.Range('D10').Value = '''0.0835'
$ws.Range('E10').Value = '  -2.07%  '
$ws.Range('D11').Value = '''32.04'
$ws.Range('E11').Value = '  -2.96%  '
$ws.Range('D13').Value = '''2.800.40'
$ws.Range('E13').Value = '  -1.82%  '
$ws.Range('E14').Value = '  -2.53%  '
$ws.Range('E15').Value = '  +0.05%  '
$ws.Range('D16').Value = '''2.440.30'
$ws.Range('E16').Value = '  -1.79%  '
$ws.Range('E17').Value = '  -1.86%  '
$ws.Range('D18').Value = '''40.999.38'
$ws.Range('E18').Value = '  -1.45%  '
$ws.Range('B19').Value = 'Uniswap'
$ws.Range('C19').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D19').Value = '''6.30'
$ws.Range('E19').Value = '  -2.48%  '
$ws.Range('B20').Value = 'ShibaInu'
$ws.Range('C20').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D20').Value = '''0.0₃0926'
$ws.Range('E20').Value = '  -2.44%  '
$ws.Range('D21').Value = '''72.29'
$ws.Range('E21').Value = '  +1.50%  '
$ws.Range('D22').Value = '''11.02'
$ws.Range('E22').Value = '  -2.80%  '
$ws.Range('D23').Value = '''234.90'
$ws.Range('E23').Value = '  -2.20%  '
$ws.Range('E24').Value = '  -1.95%  '
$ws.Range('E25').Value = '  +0.05%  '
$ws.Range('E26').Value = '  -2.99%  '
$ws.Range('D27').Value = '''24.12'
$ws.Range('E27').Value = '  -2.18%  '
$ws.Range('E28').Value = '  -1.92%  '
$ws.Range('E29').Value = '  -2.52%  '
$ws.Range('D30').Value = '''34.56'
$ws.Range('E30').Value = '  -4.54%  '
$ws.Range('D31').Value = '''158.37'
$ws.Range('E31').Value = '  -1.78%  '
$ws.Range('E32').Value = '  +0.10%  '
$ws.Range('E33').Value = '  -4.62%  '
$ws.Range('D35').Value = '''2.47'
$ws.Range('E35').Value = '  -4.38%  '
$ws.Range('D36').Value = '''16.96'
$ws.Range('E36').Value = '  -1.88%  '
$ws.Range('D37').Value = '''2.95'
$ws.Range('E37').Value = '  +1.33%  '
$ws.Range('E38').Value = '  -1.56%  '
$ws.Range('D39').Value = '''1.78'
$ws.Range('E39').Value = '  -3.83%  '
$ws.Range('D40').Value = '''0.100'
$ws.Range('E40').Value = '  -2.88%  '
$ws.Range('D42').Value = '''2.34'
$ws.Range('E42').Value = '  -4.96%  '
$ws.Range('D43').Value = '''1.995.13'
$ws.Range('E43').Value = '  +0.35%  '
$ws.Range('D44').Value = '''18.64'
$ws.Range('E44').Value = '  -1.49%  '
$ws.Range('E45').Value = '  -3.13%  '
$ws.Range('E46').Value = '  -2.86%  '
$ws.Range('D47').Value = '''9.52'
$ws.Range('E47').Value = '  +3.28%  '
$ws.Range('D48').Value = '''2.664.50'
$ws.Range('E48').Value = '  -1.60%  '
$ws.Range('D49').Value = '''94.68'
$ws.Range('E49').Value = '  -2.95%  '
$ws.Range('D50').Value = '''73.44'
$ws.Range('E50').Value = '  -1.05%  '
$ws.Range('D51').Value = '''51.95'
$ws.Range('E51').Value = '  -0.67%  '
